$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Modify existing cells per diff
$ws.Cells.Item(43, 15).Value = 2   # O43: 0 -> 2
$ws.Cells.Item(45, 18).Value = 0   # R45: blank -> 0
$ws.Cells.Item(46, 18).Value = 0   # R46: blank -> 0

# Date format used by column A (copy from existing row so new rows match)
$dateFmt = $ws.Cells.Item(46, 1).NumberFormat()

# Append new rows 47-71
# Row 47
$ws.Cells.Item(47, 1).Value = 45474
$ws.Cells.Item(47, 2).Value = 358.1499938964844
$ws.Cells.Item(47, 3).Value = 358.3999938964844
$ws.Cells.Item(47, 4).Value = 349.1000061035156
$ws.Cells.Item(47, 5).Value = 352.75
$ws.Cells.Item(47, 7).Value = 67395780
$ws.Cells.Item(47, 8).Value = 2024
$ws.Cells.Item(47, 9).Value = 7
$ws.Cells.Item(47, 10).Value = 1
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 0
$ws.Cells.Item(47, 14).Value = 27
$ws.Cells.Item(47, 15).Value = 0
$ws.Cells.Item(47, 16).Value = 0
$ws.Cells.Item(47, 17).Value = 0
$ws.Cells.Item(47, 1).NumberFormat = $dateFmt

# Row 48
$ws.Cells.Item(48, 1).Value = 45481
$ws.Cells.Item(48, 2).Value = 353
$ws.Cells.Item(48, 3).Value = 356.3999938964844
$ws.Cells.Item(48, 4).Value = 344.1000061035156
$ws.Cells.Item(48, 5).Value = 350.3500061035156
$ws.Cells.Item(48, 7).Value = 56364679
$ws.Cells.Item(48, 8).Value = 2024
$ws.Cells.Item(48, 9).Value = 7
$ws.Cells.Item(48, 10).Value = 8
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = 0
$ws.Cells.Item(48, 14).Value = 28
$ws.Cells.Item(48, 15).Value = 0
$ws.Cells.Item(48, 16).Value = 0
$ws.Cells.Item(48, 17).Value = 0
$ws.Cells.Item(48, 1).NumberFormat = $dateFmt

# Row 49
$ws.Cells.Item(49, 1).Value = 45488
$ws.Cells.Item(49, 2).Value = 354.2999877929688
$ws.Cells.Item(49, 3).Value = 356.7000122070312
$ws.Cells.Item(49, 4).Value = 332.0499877929688
$ws.Cells.Item(49, 5).Value = 338.75
$ws.Cells.Item(49, 7).Value = 81031870
$ws.Cells.Item(49, 8).Value = 2024
$ws.Cells.Item(49, 9).Value = 7
$ws.Cells.Item(49, 10).Value = 15
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = 0
$ws.Cells.Item(49, 14).Value = 29
$ws.Cells.Item(49, 15).Value = 0
$ws.Cells.Item(49, 16).Value = 0
$ws.Cells.Item(49, 17).Value = 0
$ws.Cells.Item(49, 1).NumberFormat = $dateFmt

# Row 50
$ws.Cells.Item(50, 1).Value = 45495
$ws.Cells.Item(50, 2).Value = 338.75
$ws.Cells.Item(50, 3).Value = 343
$ws.Cells.Item(50, 4).Value = 326
$ws.Cells.Item(50, 5).Value = 331.8999938964844
$ws.Cells.Item(50, 7).Value = 56761585
$ws.Cells.Item(50, 8).Value = 2024
$ws.Cells.Item(50, 9).Value = 7
$ws.Cells.Item(50, 10).Value = 22
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(50, 14).Value = 30
$ws.Cells.Item(50, 15).Value = 0
$ws.Cells.Item(50, 16).Value = 0
$ws.Cells.Item(50, 17).Value = 0
$ws.Cells.Item(50, 1).NumberFormat = $dateFmt

# Row 51
$ws.Cells.Item(51, 1).Value = 45502
$ws.Cells.Item(51, 2).Value = 334.5
$ws.Cells.Item(51, 3).Value = 335.2000122070312
$ws.Cells.Item(51, 4).Value = 324.5499877929688
$ws.Cells.Item(51, 5).Value = 326.25
$ws.Cells.Item(51, 7).Value = 54157860
$ws.Cells.Item(51, 8).Value = 2024
$ws.Cells.Item(51, 9).Value = 7
$ws.Cells.Item(51, 10).Value = 29
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(51, 14).Value = 31
$ws.Cells.Item(51, 15).Value = 0
$ws.Cells.Item(51, 16).Value = 0
$ws.Cells.Item(51, 17).Value = 0
$ws.Cells.Item(51, 1).NumberFormat = $dateFmt

# Row 52
$ws.Cells.Item(52, 1).Value = 45509
$ws.Cells.Item(52, 2).Value = 317.9500122070312
$ws.Cells.Item(52, 3).Value = 337.2000122070312
$ws.Cells.Item(52, 4).Value = 310
$ws.Cells.Item(52, 5).Value = 325.6000061035156
$ws.Cells.Item(52, 7).Value = 67590255
$ws.Cells.Item(52, 8).Value = 2024
$ws.Cells.Item(52, 9).Value = 8
$ws.Cells.Item(52, 10).Value = 5
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 13).Value = 0
$ws.Cells.Item(52, 14).Value = 32
$ws.Cells.Item(52, 15).Value = 2
$ws.Cells.Item(52, 16).Value = 0
$ws.Cells.Item(52, 17).Value = 0
$ws.Cells.Item(52, 1).NumberFormat = $dateFmt

# Row 53
$ws.Cells.Item(53, 1).Value = 45516
$ws.Cells.Item(53, 2).Value = 325
$ws.Cells.Item(53, 3).Value = 336.5
$ws.Cells.Item(53, 4).Value = 319.2000122070312
$ws.Cells.Item(53, 5).Value = 328.1499938964844
$ws.Cells.Item(53, 7).Value = 52122369
$ws.Cells.Item(53, 8).Value = 2024
$ws.Cells.Item(53, 9).Value = 8
$ws.Cells.Item(53, 10).Value = 12
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(53, 14).Value = 33
$ws.Cells.Item(53, 15).Value = 0
$ws.Cells.Item(53, 16).Value = 0
$ws.Cells.Item(53, 17).Value = 0
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt

# Row 54
$ws.Cells.Item(54, 1).Value = 45523
$ws.Cells.Item(54, 2).Value = 331.9500122070312
$ws.Cells.Item(54, 3).Value = 337.9500122070312
$ws.Cells.Item(54, 4).Value = 326.7000122070312
$ws.Cells.Item(54, 5).Value = 327.1499938964844
$ws.Cells.Item(54, 7).Value = 65701406
$ws.Cells.Item(54, 8).Value = 2024
$ws.Cells.Item(54, 9).Value = 8
$ws.Cells.Item(54, 10).Value = 19
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).Value = 0
$ws.Cells.Item(54, 14).Value = 34
$ws.Cells.Item(54, 15).Value = 0
$ws.Cells.Item(54, 16).Value = 0
$ws.Cells.Item(54, 17).Value = 0
$ws.Cells.Item(54, 1).NumberFormat = $dateFmt

# Row 55
$ws.Cells.Item(55, 1).Value = 45530
$ws.Cells.Item(55, 2).Value = 332.7000122070312
$ws.Cells.Item(55, 3).Value = 334.2000122070312
$ws.Cells.Item(55, 4).Value = 320.6499938964844
$ws.Cells.Item(55, 5).Value = 321.7000122070312
$ws.Cells.Item(55, 7).Value = 108895040
$ws.Cells.Item(55, 8).Value = 2024
$ws.Cells.Item(55, 9).Value = 8
$ws.Cells.Item(55, 10).Value = 26
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = 0
$ws.Cells.Item(55, 14).Value = 35
$ws.Cells.Item(55, 15).Value = 0
$ws.Cells.Item(55, 16).Value = 0
$ws.Cells.Item(55, 17).Value = 0
$ws.Cells.Item(55, 1).NumberFormat = $dateFmt

# Row 56
$ws.Cells.Item(56, 1).Value = 45537
$ws.Cells.Item(56, 2).Value = 323.8999938964844
$ws.Cells.Item(56, 3).Value = 355
$ws.Cells.Item(56, 4).Value = 323.1499938964844
$ws.Cells.Item(56, 5).Value = 336.8500061035156
$ws.Cells.Item(56, 7).Value = 186305220
$ws.Cells.Item(56, 8).Value = 2024
$ws.Cells.Item(56, 9).Value = 9
$ws.Cells.Item(56, 10).Value = 2
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = 0
$ws.Cells.Item(56, 14).Value = 36
$ws.Cells.Item(56, 15).Value = 0
$ws.Cells.Item(56, 16).Value = 0
$ws.Cells.Item(56, 17).Value = 0
$ws.Cells.Item(56, 1).NumberFormat = $dateFmt

# Row 57
$ws.Cells.Item(57, 1).Value = 45544
$ws.Cells.Item(57, 2).Value = 337.9500122070312
$ws.Cells.Item(57, 3).Value = 359.75
$ws.Cells.Item(57, 4).Value = 329.2999877929688
$ws.Cells.Item(57, 5).Value = 351.8999938964844
$ws.Cells.Item(57, 7).Value = 138900190
$ws.Cells.Item(57, 8).Value = 2024
$ws.Cells.Item(57, 9).Value = 9
$ws.Cells.Item(57, 10).Value = 9
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = 0
$ws.Cells.Item(57, 14).Value = 37
$ws.Cells.Item(57, 15).Value = 0
$ws.Cells.Item(57, 16).Value = 0
$ws.Cells.Item(57, 17).Value = 0
$ws.Cells.Item(57, 1).NumberFormat = $dateFmt

# Row 58
$ws.Cells.Item(58, 1).Value = 45551
$ws.Cells.Item(58, 2).Value = 352
$ws.Cells.Item(58, 3).Value = 355.8999938964844
$ws.Cells.Item(58, 4).Value = 338.4500122070312
$ws.Cells.Item(58, 5).Value = 353.8999938964844
$ws.Cells.Item(58, 7).Value = 94499190
$ws.Cells.Item(58, 8).Value = 2024
$ws.Cells.Item(58, 9).Value = 9
$ws.Cells.Item(58, 10).Value = 16
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = 0
$ws.Cells.Item(58, 14).Value = 38
$ws.Cells.Item(58, 15).Value = 0
$ws.Cells.Item(58, 16).Value = 1
$ws.Cells.Item(58, 17).Value = 1
$ws.Cells.Item(58, 1).NumberFormat = $dateFmt

# Row 59
$ws.Cells.Item(59, 1).Value = 45558
$ws.Cells.Item(59, 2).Value = 354.3500061035156
$ws.Cells.Item(59, 3).Value = 363
$ws.Cells.Item(59, 4).Value = 348.9500122070312
$ws.Cells.Item(59, 5).Value = 359.1499938964844
$ws.Cells.Item(59, 7).Value = 129753300
$ws.Cells.Item(59, 8).Value = 2024
$ws.Cells.Item(59, 9).Value = 9
$ws.Cells.Item(59, 10).Value = 23
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 13).Value = 0
$ws.Cells.Item(59, 14).Value = 39
$ws.Cells.Item(59, 15).Value = 1
$ws.Cells.Item(59, 16).Value = 0
$ws.Cells.Item(59, 17).Value = 0
$ws.Cells.Item(59, 1).NumberFormat = $dateFmt

# Row 60
$ws.Cells.Item(60, 1).Value = 45565
$ws.Cells.Item(60, 2).Value = 357
$ws.Cells.Item(60, 3).Value = 358.1499938964844
$ws.Cells.Item(60, 4).Value = 336.5
$ws.Cells.Item(60, 5).Value = 338.7999877929688
$ws.Cells.Item(60, 7).Value = 86363250
$ws.Cells.Item(60, 8).Value = 2024
$ws.Cells.Item(60, 9).Value = 9
$ws.Cells.Item(60, 10).Value = 30
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = 0
$ws.Cells.Item(60, 14).Value = 40
$ws.Cells.Item(60, 15).Value = 0
$ws.Cells.Item(60, 16).Value = 0
$ws.Cells.Item(60, 17).Value = 0
$ws.Cells.Item(60, 1).NumberFormat = $dateFmt

# Row 61
$ws.Cells.Item(61, 1).Value = 45572
$ws.Cells.Item(61, 2).Value = 346.1000061035156
$ws.Cells.Item(61, 3).Value = 349.7000122070312
$ws.Cells.Item(61, 4).Value = 329.2000122070312
$ws.Cells.Item(61, 5).Value = 341.5
$ws.Cells.Item(61, 7).Value = 88568027
$ws.Cells.Item(61, 8).Value = 2024
$ws.Cells.Item(61, 9).Value = 10
$ws.Cells.Item(61, 10).Value = 7
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = 0
$ws.Cells.Item(61, 14).Value = 41
$ws.Cells.Item(61, 15).Value = 0
$ws.Cells.Item(61, 16).Value = 0
$ws.Cells.Item(61, 17).Value = 0
$ws.Cells.Item(61, 1).NumberFormat = $dateFmt

# Row 62
$ws.Cells.Item(62, 1).Value = 45579
$ws.Cells.Item(62, 2).Value = 342.7999877929688
$ws.Cells.Item(62, 3).Value = 343.4500122070312
$ws.Cells.Item(62, 4).Value = 325.5
$ws.Cells.Item(62, 5).Value = 330.3500061035156
$ws.Cells.Item(62, 7).Value = 46746485
$ws.Cells.Item(62, 8).Value = 2024
$ws.Cells.Item(62, 9).Value = 10
$ws.Cells.Item(62, 10).Value = 14
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = 0
$ws.Cells.Item(62, 14).Value = 42
$ws.Cells.Item(62, 15).Value = 0
$ws.Cells.Item(62, 16).Value = 0
$ws.Cells.Item(62, 17).Value = 0
$ws.Cells.Item(62, 1).NumberFormat = $dateFmt

# Row 63
$ws.Cells.Item(63, 1).Value = 45586
$ws.Cells.Item(63, 2).Value = 334.3999938964844
$ws.Cells.Item(63, 3).Value = 334.3999938964844
$ws.Cells.Item(63, 4).Value = 306
$ws.Cells.Item(63, 5).Value = 311.2000122070312
$ws.Cells.Item(63, 7).Value = 60601997
$ws.Cells.Item(63, 8).Value = 2024
$ws.Cells.Item(63, 9).Value = 10
$ws.Cells.Item(63, 10).Value = 21
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(63, 14).Value = 43
$ws.Cells.Item(63, 15).Value = 0
$ws.Cells.Item(63, 16).Value = 0
$ws.Cells.Item(63, 17).Value = 0
$ws.Cells.Item(63, 1).NumberFormat = $dateFmt

# Row 64
$ws.Cells.Item(64, 1).Value = 45593
$ws.Cells.Item(64, 2).Value = 314.5
$ws.Cells.Item(64, 3).Value = 326.9500122070312
$ws.Cells.Item(64, 4).Value = 310.6499938964844
$ws.Cells.Item(64, 5).Value = 326.2000122070312
$ws.Cells.Item(64, 7).Value = 35227255
$ws.Cells.Item(64, 8).Value = 2024
$ws.Cells.Item(64, 9).Value = 10
$ws.Cells.Item(64, 10).Value = 28
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(64, 14).Value = 44
$ws.Cells.Item(64, 15).Value = 0
$ws.Cells.Item(64, 16).Value = 0
$ws.Cells.Item(64, 17).Value = 0
$ws.Cells.Item(64, 1).NumberFormat = $dateFmt

# Row 65
$ws.Cells.Item(65, 1).Value = 45600
$ws.Cells.Item(65, 2).Value = 326.7999877929688
$ws.Cells.Item(65, 3).Value = 328.5
$ws.Cells.Item(65, 4).Value = 314.7999877929688
$ws.Cells.Item(65, 5).Value = 316.0499877929688
$ws.Cells.Item(65, 7).Value = 30218929
$ws.Cells.Item(65, 8).Value = 2024
$ws.Cells.Item(65, 9).Value = 11
$ws.Cells.Item(65, 10).Value = 4
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(65, 14).Value = 45
$ws.Cells.Item(65, 15).Value = 0
$ws.Cells.Item(65, 16).Value = 0
$ws.Cells.Item(65, 17).Value = 2
$ws.Cells.Item(65, 1).NumberFormat = $dateFmt

# Row 66
$ws.Cells.Item(66, 1).Value = 45607
$ws.Cells.Item(66, 2).Value = 314
$ws.Cells.Item(66, 3).Value = 320.5
$ws.Cells.Item(66, 4).Value = 298
$ws.Cells.Item(66, 5).Value = 318.3500061035156
$ws.Cells.Item(66, 7).Value = 71508852
$ws.Cells.Item(66, 8).Value = 2024
$ws.Cells.Item(66, 9).Value = 11
$ws.Cells.Item(66, 10).Value = 11
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = 0
$ws.Cells.Item(66, 14).Value = 46
$ws.Cells.Item(66, 15).Value = 2
$ws.Cells.Item(66, 16).Value = 0
$ws.Cells.Item(66, 17).Value = 0
$ws.Cells.Item(66, 1).NumberFormat = $dateFmt

# Row 67
$ws.Cells.Item(67, 1).Value = 45614
$ws.Cells.Item(67, 2).Value = 320
$ws.Cells.Item(67, 3).Value = 323.7999877929688
$ws.Cells.Item(67, 4).Value = 312.3500061035156
$ws.Cells.Item(67, 5).Value = 313.7999877929688
$ws.Cells.Item(67, 7).Value = 46322349
$ws.Cells.Item(67, 8).Value = 2024
$ws.Cells.Item(67, 9).Value = 11
$ws.Cells.Item(67, 10).Value = 18
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = 0
$ws.Cells.Item(67, 14).Value = 47
$ws.Cells.Item(67, 15).Value = 0
$ws.Cells.Item(67, 16).Value = 0
$ws.Cells.Item(67, 17).Value = 0
$ws.Cells.Item(67, 1).NumberFormat = $dateFmt

# Row 68
$ws.Cells.Item(68, 1).Value = 45621
$ws.Cells.Item(68, 2).Value = 318.8500061035156
$ws.Cells.Item(68, 3).Value = 332.2999877929688
$ws.Cells.Item(68, 4).Value = 318
$ws.Cells.Item(68, 5).Value = 328.3500061035156
$ws.Cells.Item(68, 7).Value = 76004136
$ws.Cells.Item(68, 8).Value = 2024
$ws.Cells.Item(68, 9).Value = 11
$ws.Cells.Item(68, 10).Value = 25
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(68, 14).Value = 48
$ws.Cells.Item(68, 15).Value = 0
$ws.Cells.Item(68, 16).Value = 0
$ws.Cells.Item(68, 17).Value = 0
$ws.Cells.Item(68, 1).NumberFormat = $dateFmt

# Row 69
$ws.Cells.Item(69, 1).Value = 45628
$ws.Cells.Item(69, 2).Value = 326
$ws.Cells.Item(69, 3).Value = 347.3500061035156
$ws.Cells.Item(69, 4).Value = 325.1000061035156
$ws.Cells.Item(69, 5).Value = 338.8999938964844
$ws.Cells.Item(69, 7).Value = 63329885
$ws.Cells.Item(69, 8).Value = 2024
$ws.Cells.Item(69, 9).Value = 12
$ws.Cells.Item(69, 10).Value = 2
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(69, 14).Value = 49
$ws.Cells.Item(69, 15).Value = 0
$ws.Cells.Item(69, 16).Value = 0
$ws.Cells.Item(69, 17).Value = 0
$ws.Cells.Item(69, 1).NumberFormat = $dateFmt

# Row 70
$ws.Cells.Item(70, 1).Value = 45635
$ws.Cells.Item(70, 2).Value = 338.5
$ws.Cells.Item(70, 3).Value = 346.2999877929688
$ws.Cells.Item(70, 4).Value = 330.6499938964844
$ws.Cells.Item(70, 5).Value = 339.75
$ws.Cells.Item(70, 7).Value = 57821189
$ws.Cells.Item(70, 8).Value = 2024
$ws.Cells.Item(70, 9).Value = 12
$ws.Cells.Item(70, 10).Value = 9
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(70, 14).Value = 50
$ws.Cells.Item(70, 15).Value = 0
$ws.Cells.Item(70, 16).Value = 0
$ws.Cells.Item(70, 17).Value = 0
$ws.Cells.Item(70, 1).NumberFormat = $dateFmt

# Row 71
$ws.Cells.Item(71, 1).Value = 45642
$ws.Cells.Item(71, 2).Value = 339.25
$ws.Cells.Item(71, 3).Value = 342.4500122070312
$ws.Cells.Item(71, 4).Value = 302.7999877929688
$ws.Cells.Item(71, 5).Value = 304.2999877929688
$ws.Cells.Item(71, 7).Value = 69768815
$ws.Cells.Item(71, 8).Value = 2024
$ws.Cells.Item(71, 9).Value = 12
$ws.Cells.Item(71, 10).Value = 16
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(71, 14).Value = 51
$ws.Cells.Item(71, 15).Value = 0
$ws.Cells.Item(71, 16).Value = 0
$ws.Cells.Item(71, 17).Value = 0
$ws.Cells.Item(71, 1).NumberFormat = $dateFmt
